# "allow player to move in 4 directions"
# The TransformComponent's "button press" row (right/left) is split so a
# second row is added for up/down, and the position-of-x/position-of-y rows
# get their own (non-blank) id values. Row numbers shift: the old rows
# 2,3,4,5 become rows 4,5,7,8,9 (row 6 stays blank), while the
# SpriteComponent rows (15-17) keep their row numbers but get new ids.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old TransformComponent block (previously rows 2-5) before
# rewriting it across its new row range (4-9).
$ws.Range("A2:E9").Clear() | Out-Null

# Row 4: button press 1 = right / -1 = left;  (id changed 123 -> 111)
$ws.Cells.Item(4, 1).Value = "TransformComponent"
$ws.Cells.Item(4, 2).Value = 111
$ws.Cells.Item(4, 3).Value = "[-1,1]"
$ws.Cells.Item(4, 4).Value = "button press 1 = right / -1 = left;"

# Row 5: NEW - button press 1 = up / -1 = down;
$ws.Cells.Item(5, 1).Value = "TransformComponent"
$ws.Cells.Item(5, 2).Value = 112
$ws.Cells.Item(5, 3).Value = "[-1,1]"
$ws.Cells.Item(5, 4).Value = "button press 1 = up / -1 = down;"

# Row 7: dt for speed in ms. (id unchanged, 124)
$ws.Cells.Item(7, 1).Value = "TransformComponent"
$ws.Cells.Item(7, 2).Value = 124
$ws.Cells.Item(7, 3).Value = "Integer"
$ws.Cells.Item(7, 4).Value = "dt for speed in ms."

# Row 8: position of x (now with an id + Integer type)
$ws.Cells.Item(8, 1).Value = "TransformComponent"
$ws.Cells.Item(8, 2).Value = 101
$ws.Cells.Item(8, 3).Value = "Integer"
$ws.Cells.Item(8, 4).Value = "position of x"
$ws.Cells.Item(8, 5).Value = $true

# Row 9: position of y (now with an id + Integer type)
$ws.Cells.Item(9, 1).Value = "TransformComponent"
$ws.Cells.Item(9, 2).Value = 102
$ws.Cells.Item(9, 3).Value = "Integer"
$ws.Cells.Item(9, 4).Value = "position of y"

# SpriteComponent id renumbering (rows keep their position, only B changes)
$ws.Cells.Item(15, 2).Value = 211
$ws.Cells.Item(16, 2).Value = 212
$ws.Cells.Item(17, 2).Value = 213

# Match the saved selection from the edit
$ws.Range("A7:C8").Select()
